# Schedule.xlsx update -- "Important critical updates - Working live with client"
#
# - Marks the last "Ongoing" task as "Done"
# - Reduces several Hours entries (re-estimated effort)
# - Adds two new rows: Revision-06 (bug fixes) and Revision-07 (previous file
#   check & bug fixes) at the bottom of the schedule table
# - Extends the Hours total formula to cover the new rows
# - Bolds the "Sr. #" column for the data rows
# - Shrinks the thin spacer row under the header
# - Resets the view (scroll position/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank rows just above the existing totals row (row 24),
#    so the totals row + the two spacer rows below it shift down by two.
# ---------------------------------------------------------------------------
$ws.Rows("24:25").Insert()

# ---------------------------------------------------------------------------
# 2. Re-key the Hours column (effort re-estimation) for the existing tasks.
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("F12").Value = 1.5
$ws.Range("F13").Value = 1.5
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 4
$ws.Range("F21").Value = 3

# The last existing task is finished now -- flip its status from "Ongoing"
# to "Done".
$ws.Range("E23").Value = "Done"

# ---------------------------------------------------------------------------
# 3. Fill in the two new schedule rows.
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = 18
$ws.Range("C24").Value = "Revision-06 - Bug fixes"
$ws.Range("D24").Value = "31 - 12 - 2019"
$ws.Range("E24").Value = "Done"
$ws.Range("F24").Value = 1

$ws.Range("B25").Value = 19
$ws.Range("C25").Value = "Revision-07 - Previous file check & bug fixes"
$ws.Range("D25").Value = "02 - 01 - 2020"
$ws.Range("E25").Value = "Done"
$ws.Range("F25").Value = 1

# ---------------------------------------------------------------------------
# 4. Extend the Hours total to include the two new rows (now on row 26).
# ---------------------------------------------------------------------------
$ws.Range("F26").Formula = "=SUM(F7:F25)"

# ---------------------------------------------------------------------------
# 5. Match formatting tweaks: bold the Sr. # column for all the data rows,
#    and shrink the thin divider row right under the table header.
# ---------------------------------------------------------------------------
$ws.Range("B7:B25").Font.Bold = $true
$ws.Rows(6).RowHeight = 9.75

# ---------------------------------------------------------------------------
# 6. Reset the view: scroll back to the top and select C16.
# ---------------------------------------------------------------------------
$ws.Range("C16").Select()
